$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in Senado 18's color (row 10, column C): #348cae4 -> #348cae
$ws.Range("C10").Value = "#348cae"

# Insert a new row at position 15 (shifts Distrito local 21 etc. down)
$ws.Rows.Item(15).Insert()

# Populate the new row 15 with the Gobernatura 21 entry
$ws.Range("A15").Value = "Gobernatura 21"
$ws.Range("B15").Value = "gb_21"
$ws.Range("C15").Value = "#4361ee"

# Update the active selection to match the authored state
$ws.Range("G14").Select()
